$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append three new daily rows (05-10-2021, 06-10-2021, 07-10-2021) after the
# existing data, which previously ended at row 191 (04-10-2021).
#
# Column A holds date-like text labels ("05-10-2021", ...). If entered as a
# plain .Value, Excel's smart-entry parsing would turn them into real dates.
# To keep them as literal text (matching the rest of the "Serie" column,
# which stores plain text labels with no special number formatting), each
# label is entered as a ="..." text formula and then converted in place to
# a static value via Copy / PasteSpecial values-only. That yields a plain
# shared-string cell with no formula left behind and no extra cell styling.

function Set-DateLabel($range, $text) {
    $range.Formula = '="' + $text + '"'
    $range.Copy()
    $range.PasteSpecial(-4163)  # xlPasteValues
}

Set-DateLabel $ws.Range("A192") "05-10-2021"
$ws.Range("B192").Value = 13530
$ws.Range("C192").Value = 20274
$ws.Range("D192").Value = -6744

Set-DateLabel $ws.Range("A193") "06-10-2021"
$ws.Range("B193").Value = 14066
$ws.Range("C193").Value = 20565
$ws.Range("D193").Value = -6500

Set-DateLabel $ws.Range("A194") "07-10-2021"
$ws.Range("B194").Value = 13786
$ws.Range("C194").Value = 20503
$ws.Range("D194").Value = -6717

$excel.CutCopyMode = $false
